$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.242386666666667
$ws.Range("H2").Value = 6.72716
$ws.Range("I2").Value = 0.04442500453715972
$ws.Range("J2").Value = 0.04442500453715972
$ws.Range("M2").Value = 42.09975866666667
$ws.Range("N2").Value = 126.299276
$ws.Range("O2").Value = 0.3315552933456474
$ws.Range("P2").Value = 0.3315552933456474
$ws.Range("Q2").Value = 94.40393750401778
$ws.Range("R2").Value = 849.6354375361599
$ws.Range("S2").Value = 0.01472934541119971
$ws.Range("T2").Value = 0.01472934541119971
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.242386666666667
$ws.Range("H3").Value = 6.72716
$ws.Range("I3").Value = 0.04442500453715972
$ws.Range("J3").Value = 0.04442500453715972
$ws.Range("O3").Value = 0.4502223747274475
$ws.Range("P3").Value = 0.4502223747274475
$ws.Range("Q3").Value = 128.1920867490756
$ws.Range("R3").Value = 1153.72878074168
$ws.Range("S3").Value = 0.02000113103999768
$ws.Range("T3").Value = 0.02000113103999768
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.242386666666667
$ws.Range("H4").Value = 6.72716
$ws.Range("I4").Value = 0.04442500453715972
$ws.Range("J4").Value = 0.04442500453715972
$ws.Range("M4").Value = 27.596267
$ws.Range("N4").Value = 82.78880100000001
$ws.Range("O4").Value = 0.2173335118824389
$ws.Range("P4").Value = 0.2173335118824389
$ws.Range("Q4").Value = 61.88150117057333
$ws.Range("R4").Value = 556.9335105351601
$ws.Range("S4").Value = 0.009655042251454202
$ws.Range("T4").Value = 0.009655042251454202
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.242386666666667
$ws.Range("H5").Value = 6.72716
$ws.Range("I5").Value = 0.04442500453715972
$ws.Range("J5").Value = 0.04442500453715972
$ws.Range("M5").Value = 0.1128593333333333
$ws.Range("N5").Value = 0.338578
$ws.Range("O5").Value = 0.0008888200444663087
$ws.Range("P5").Value = 0.0008888200444663087
$ws.Range("Q5").Value = 0.2530742642755555
$ws.Range("R5").Value = 2.27766837848
$ws.Range("S5").Value = 0.00003948583450813427
$ws.Range("T5").Value = 0.00003948583450813427
$ws.Range("I6").Value = 0.4052409520727612
$ws.Range("J6").Value = 0.4052409520727612
$ws.Range("M6").Value = 42.09975866666667
$ws.Range("N6").Value = 126.299276
$ws.Range("O6").Value = 0.3315552933456474
$ws.Range("P6").Value = 0.3315552933456474
$ws.Range("Q6").Value = 861.1443467956368
$ws.Range("R6").Value = 7750.299121160731
$ws.Range("S6").Value = 0.1343597827401538
$ws.Range("T6").Value = 0.1343597827401538
$ws.Range("I7").Value = 0.4052409520727612
$ws.Range("J7").Value = 0.4052409520727612
$ws.Range("O7").Value = 0.4502223747274475
$ws.Range("P7").Value = 0.4502223747274475
$ws.Range("S7").Value = 0.1824485437790103
$ws.Range("T7").Value = 0.1824485437790103
$ws.Range("I8").Value = 0.4052409520727612
$ws.Range("J8").Value = 0.4052409520727612
$ws.Range("M8").Value = 27.596267
$ws.Range("N8").Value = 82.78880100000001
$ws.Range("O8").Value = 0.2173335118824389
$ws.Range("P8").Value = 0.2173335118824389
$ws.Range("Q8").Value = 564.4775664362396
$ws.Range("R8").Value = 5080.298097926157
$ws.Range("S8").Value = 0.08807243927255629
$ws.Range("T8").Value = 0.08807243927255629
$ws.Range("I9").Value = 0.4052409520727612
$ws.Range("J9").Value = 0.4052409520727612
$ws.Range("M9").Value = 0.1128593333333333
$ws.Range("N9").Value = 0.338578
$ws.Range("O9").Value = 0.0008888200444663087
$ws.Range("P9").Value = 0.0008888200444663087
$ws.Range("Q9").Value = 2.308520997771777
$ws.Range("R9").Value = 20.776688979946
$ws.Range("S9").Value = 0.0003601862810408809
$ws.Range("T9").Value = 0.0003601862810408809
$ws.Range("G10").Value = 27.778539
$ws.Range("H10").Value = 83.335617
$ws.Range("I10").Value = 0.5503340433900792
$ws.Range("J10").Value = 0.5503340433900791
$ws.Range("M10").Value = 42.09975866666667
$ws.Range("N10").Value = 126.299276
$ws.Range("O10").Value = 0.3315552933456474
$ws.Range("P10").Value = 0.3315552933456474
$ws.Range("Q10").Value = 1169.469788012588
$ws.Range("R10").Value = 10525.22809211329
$ws.Range("S10").Value = 0.1824661651942939
$ws.Range("T10").Value = 0.1824661651942939
$ws.Range("G11").Value = 27.778539
$ws.Range("H11").Value = 83.335617
$ws.Range("I11").Value = 0.5503340433900792
$ws.Range("J11").Value = 0.5503340433900791
$ws.Range("O11").Value = 0.4502223747274475
$ws.Range("P11").Value = 0.4502223747274475
$ws.Range("Q11").Value = 1588.035165471274
$ws.Range("R11").Value = 14292.31648924147
$ws.Range("S11").Value = 0.2477726999084396
$ws.Range("T11").Value = 0.2477726999084396
$ws.Range("G12").Value = 27.778539
$ws.Range("H12").Value = 83.335617
$ws.Range("I12").Value = 0.5503340433900792
$ws.Range("J12").Value = 0.5503340433900791
$ws.Range("M12").Value = 27.596267
$ws.Range("N12").Value = 82.78880100000001
$ws.Range("O12").Value = 0.2173335118824389
$ws.Range("P12").Value = 0.2173335118824389
$ws.Range("Q12").Value = 766.583979113913
$ws.Range("R12").Value = 6899.255812025218
$ws.Range("S12").Value = 0.1196060303584284
$ws.Range("T12").Value = 0.1196060303584284
$ws.Range("G13").Value = 27.778539
$ws.Range("H13").Value = 83.335617
$ws.Range("I13").Value = 0.5503340433900792
$ws.Range("J13").Value = 0.5503340433900791
$ws.Range("M13").Value = 0.1128593333333333
$ws.Range("N13").Value = 0.338578
$ws.Range("O13").Value = 0.0008888200444663087
$ws.Range("P13").Value = 0.0008888200444663087
$ws.Range("Q13").Value = 3.135067392514
$ws.Range("R13").Value = 28.215606532626
$ws.Range("S13").Value = 0.0004891479289172937
$ws.Range("T13").Value = 0.0004891479289172936
